# AfDD_2022_Annex_Table_Tab10.xlsx refresh
# - Bumps the window size recorded for the workbook view.
# - Updates the "Source:" footnote with the newer ILO/WB retrieval dates.
# - Refreshes Oil/Mineral/Total natural-resources rents (% of GDP) figures
#   (columns C, D, E) for the regional & grouping aggregate rows, reflecting
#   upstream data revisions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab10")

# --- Workbook window geometry -------------------------------------------------
$win = $excel.ActiveWindow
$win.Width = 19200
$win.Height = 10400

# --- Footnote / source text ---------------------------------------------------
$ws.Range("A105").Value = "Source: International Labour Organisation (retrieved 14/11/2021), United Nations Statistics Division National Accounts (Analysis of Main Aggregates dataset uploaded in December 2020), World Bank World Development Indicators (database and data releases from central banks, national statistical agencies, and World Bank country desks -22/04/2022)."

# --- Updated rents data (Oil / Mineral / Total natural resources rents, % GDP) --

$ws.Range("C13").Value = 4.24511870492843
$ws.Range("D13").Value = 1.02561809544274
$ws.Range("E13").Value = 7.87601762055452
$ws.Range("C23").Value = 9.22394085012542
$ws.Range("D23").Value = 0.93080563470538
$ws.Range("E23").Value = 15.0418374652097
$ws.Range("C38").Value = 0.64473634843504
$ws.Range("D38").Value = 0.00052466462392
$ws.Range("E38").Value = 3.48985318086668
$ws.Range("C45").Value = 7.32151739129707
$ws.Range("D45").Value = 0.15644886761608
$ws.Range("E45").Value = 8.53114524107064
$ws.Range("C61").Value = 5.2818945746942
$ws.Range("D61").Value = 0.07097767178881
$ws.Range("E61").Value = 7.48652374662629
$ws.Range("C62").Value = 5.22168366650812
$ws.Range("D62").Value = 0.31636727112377
$ws.Range("E62").Value = 7.6403140161356
$ws.Range("C63").Value = 1.53799932133146
$ws.Range("D63").Value = 0.21289374191068
$ws.Range("E63").Value = 2.39689202036949
$ws.Range("C64").Value = 1.77368365253415
$ws.Range("D64").Value = 0.44582365183316
$ws.Range("E64").Value = 2.75337319723518
$ws.Range("C65").Value = 0.89492303171264
$ws.Range("D65").Value = 0.30388181616974
$ws.Range("E65").Value = 2.15461165739714
$ws.Range("C66").Value = 1.725631105694
$ws.Range("D66").Value = 0.21816406658472
$ws.Range("E66").Value = 2.66396067581068
$ws.Range("C67").Value = 3.86105380973162
$ws.Range("D67").Value = 0.16784557187966
$ws.Range("E67").Value = 6.07552477629987
$ws.Range("C68").Value = 5.278105327627
$ws.Range("D68").Value = 0.12880948769287
$ws.Range("E68").Value = 6.88485843356153
$ws.Range("C69").Value = 0.00981467579527
$ws.Range("D69").Value = 0.00004643947943
$ws.Range("E69").Value = 2.57205177967192
$ws.Range("C70").Value = 14.8686072250135
$ws.Range("D70").Value = 0.52435270176506
$ws.Range("E70").Value = 18.7344191769131
$ws.Range("C71").Value = 5.2818945746942
$ws.Range("D71").Value = 0.07097767178881
$ws.Range("E71").Value = 7.48652374662629
$ws.Range("C72").Value = 0.85322440258808
$ws.Range("D72").Value = 0.00070003148633
$ws.Range("E72").Value = 3.92029336585616
$ws.Range("C73").Value = 3.416424340786
$ws.Range("D73").Value = 0.99593010995525
$ws.Range("E73").Value = 7.20485373593153
$ws.Range("C74").Value = 11.1257606541548
$ws.Range("D74").Value = 0.3293511332772
$ws.Range("E74").Value = 12.6231677835469
$ws.Range("C75").Value = 20.6468361476981
$ws.Range("E75").Value = 23.8896323141622
$ws.Range("C76").Value = 0.90729207139557
$ws.Range("D76").Value = 0.1737090100106
$ws.Range("E76").Value = 2.78864700561867
$ws.Range("C77").Value = 1.98187015905833
$ws.Range("D77").Value = 0.58838826316593
$ws.Range("E77").Value = 3.23041968769066
$ws.Range("C78").Value = 0.04016033223356
$ws.Range("D78").Value = 0.03012799355049
$ws.Range("E78").Value = 0.19059421440465
$ws.Range("C79").Value = 0.36314711207219
$ws.Range("D79").Value = 0.11152221306819
$ws.Range("E79").Value = 0.67860633949165
$ws.Range("C80").Value = 13.5989758264585
$ws.Range("D80").Value = 0.15081776776983
$ws.Range("E80").Value = 15.6061513885532
$ws.Range("C81").Value = 15.5253439794699
$ws.Range("D81").Value = 0.33605502046029
$ws.Range("E81").Value = 18.3095829784258
$ws.Range("C82").Value = 1.55637400376145
$ws.Range("D82").Value = 0.3888000074687
$ws.Range("E82").Value = 4.15502782394869
$ws.Range("C83").Value = 0.45784838059269
$ws.Range("D83").Value = 0.20338319870721
$ws.Range("E83").Value = 1.16811024234402
$ws.Range("C84").Value = 1.22011649852953
$ws.Range("D84").Value = 0.389526773682
$ws.Range("E84").Value = 6.18806120123897
$ws.Range("C86").Value = 6.19966599264777
$ws.Range("D86").Value = 0.11022997292483
$ws.Range("E86").Value = 7.79523958863848
$ws.Range("C87").Value = 0.49907381428577
$ws.Range("D87").Value = 0.43433652213344
$ws.Range("E87").Value = 2.18584924464007
$ws.Range("C88").Value = 4.99783018981369
$ws.Range("D88").Value = 1.1120113589761
$ws.Range("E88").Value = 8.41436165662549
$ws.Range("C89").Value = 2.20568353934627
$ws.Range("D89").Value = 0.27170604001456
$ws.Range("E89").Value = 3.42722125236528
$ws.Range("C90").Value = 1.37953764164502
$ws.Range("D90").Value = 0.10298128054202
$ws.Range("E90").Value = 1.72443460230978
$ws.Range("C91").Value = 4.17778920081595
$ws.Range("D91").Value = 0.5023067614396
$ws.Range("E91").Value = 8.2815665605199
$ws.Range("C92").Value = 0.30778174802583
$ws.Range("D92").Value = 0.14060517761255
$ws.Range("E92").Value = 1.65738785692054
$ws.Range("E93").Value = 0.9664676783884
$ws.Range("C94").Value = 0.47485830524745
$ws.Range("D94").Value = 0.09068611735378
$ws.Range("E94").Value = 1.31197037483654
$ws.Range("C95").Value = 0.6847124534221
$ws.Range("D95").Value = 0.19222733364474
$ws.Range("E95").Value = 5.00904502095678
$ws.Range("C96").Value = 6.84549859618508
$ws.Range("D96").Value = 1.04244520231506
$ws.Range("E96").Value = 10.461663234768
$ws.Range("C97").Value = 6.49390830861809
$ws.Range("D97").Value = 0.26345531582738
$ws.Range("E97").Value = 9.43270101135141
$ws.Range("C98").Value = 5.69821559519683
$ws.Range("D98").Value = 0.11983896057817
$ws.Range("E98").Value = 6.70516652327677
$ws.Range("C99").Value = 21.0021688369151
$ws.Range("D99").Value = 0.30032705968755
$ws.Range("E99").Value = 22.9066097067589
